$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "season record" columns appended after the existing data (A:AC).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold / centered / bordered header formatting used by A1:AC1.
$headerRng = $ws.Range("AD1:AF1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# Every player row (2-50) gets the team's 2001 season record.
$ws.Range("AD2:AD50").Value = 82
$ws.Range("AE2:AE50").Value = 79
$ws.Range("AF2:AF50").Value = 0
